$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Formula = '=IF(ISERROR(B3),"ERROR",IF(ISERROR(C3),"FAIL",IF(B3=C3,"PASS","FAIL")))'
$ws.Range("D4:D16").Formula = '=IF(ISERROR(B4),"ERROR",IF(ISERROR(C4),"FAIL",IF(B4=C4,"PASS","FAIL")))'
$ws.Range("D16").Style = $ws.Range("D15").Style
